$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# The Date column (H2:H10) stores its dates as plain text, so make sure
# Excel keeps writing text (not auto-converted serial date numbers).
$dates = $ws.Range("H2:H10")
$originalFormat = $dates.NumberFormat
$dates.NumberFormat = "@"

$ws.Range("H2").Value = "1.5.2025"
$ws.Range("H3").Value = "2.5.2025"
$ws.Range("H4").Value = "30.4.2025"
$ws.Range("H5").Value = "2.5.2025"
$ws.Range("H6").Value = "4.5.2025"
$ws.Range("H7").Value = "2.5.2025"
$ws.Range("H8").Value = "3.5.2025"
$ws.Range("H9").Value = "4.5.2025"
$ws.Range("H10").Value = "2.5.2025"

# Restore the original cell number format (dd.mm.yyyy) so the saved
# style indexes match the source workbook.
$dates.NumberFormat = $originalFormat

# Move the active selection to H10 as in the saved file
$ws.Range("H10").Select()
